$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (English record): Vostro / Dell / DKS laptop spec
$ws.Range("B2").Value = "Vostro"
$ws.Range("C2").Value = "Dell"
$ws.Range("D2").Value = 3568
$ws.Range("E2").Value = "DKS"
$ws.Range("F2").Value = 1.454
$ws.Range("G2").Value = "To take enrollments"
$ws.Range("H2").Value = "eng"

# Update row 3 (Arabic record): same machine, translated fields
$ws.Range("D3").Value = 3568
$ws.Range("E3").Value = "DKS"
$ws.Range("F3").Value = 1.454
$ws.Range("H3").Value = "ara"

# Adjust the view: scroll so column C is the left-most visible column,
# and select from A4 down to the bottom of the sheet
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("A4:XFD1048576").Select
$ws.Range("C4").Activate

# Configure page setup for printing
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
